$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Characters(21, 2).Text = "35"
$ws.Range("C9").Characters(27, 9).Text = "8/28/2023"
$ws.Range("C9").Characters(47, 9).Text = "9/3/2023"

# --- Plain numeric value updates ---
$ws.Range("M15").Value = 63.636363636363
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 178
$ws.Range("J16").Value = 135
$ws.Range("K16").Value = 31.851851851851
$ws.Range("L16").Value = 128.205128205128
$ws.Range("M16").Value = 36.923076923076
$ws.Range("N16").Value = -74.607703281027
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 700
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 160
$ws.Range("J17").Value = 152
$ws.Range("K17").Value = 5.263157894736
$ws.Range("L17").Value = 16.788321167883
$ws.Range("M17").Value = 81.818181818181
$ws.Range("N17").Value = -24.882629107981
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 7.142857142857
$ws.Range("J18").Value = 117
$ws.Range("K18").Value = 13.675213675213
$ws.Range("L18").Value = 12.711864406779
$ws.Range("M18").Value = -17.901234567901
$ws.Range("N18").Value = -86.686686686686
$ws.Range("C19").Value = 20
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = -8.695652173913
$ws.Range("I19").Value = 507
$ws.Range("J19").Value = 469
$ws.Range("K19").Value = 8.102345415778
$ws.Range("L19").Value = 81.071428571428
$ws.Range("M19").Value = 67.326732673267
$ws.Range("N19").Value = -18.488745980707
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 33
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 37.5
$ws.Range("I20").Value = 180
$ws.Range("J20").Value = 138
$ws.Range("K20").Value = 30.434782608695
$ws.Range("L20").Value = 45.16129032258
$ws.Range("M20").Value = 21.621621621621
$ws.Range("N20").Value = -86.666666666666
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 9.375
$ws.Range("F21").Value = 150
$ws.Range("H21").Value = 2.739726027397
$ws.Range("I21").Value = 1176
$ws.Range("J21").Value = 1026
$ws.Range("K21").Value = 14.619883040935
$ws.Range("L21").Value = 57.429718875502
$ws.Range("M21").Value = 39.336492890995
$ws.Range("N21").Value = -69.846153846153
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = -45.454545454545
$ws.Range("I22").Value = 61
$ws.Range("J22").Value = 55
$ws.Range("K22").Value = 10.90909090909
$ws.Range("L22").Value = 134.615384615385
$ws.Range("M22").Value = 96.666666666666
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -5.882352941176
$ws.Range("F24").Value = 173
$ws.Range("G24").Value = 152
$ws.Range("H24").Value = 13.815789473684
$ws.Range("I24").Value = 1334
$ws.Range("J24").Value = 965
$ws.Range("K24").Value = 38.238341968911
$ws.Range("L24").Value = 61.305925030229
$ws.Range("M24").Value = 105.546995377504
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 36
$ws.Range("H25").Value = 5.882352941176
$ws.Range("I25").Value = 342
$ws.Range("J25").Value = 384
$ws.Range("K25").Value = -10.9375
$ws.Range("L25").Value = 8.571428571428
$ws.Range("M25").Value = -1.156069364161
$ws.Range("I26").Value = 24
$ws.Range("K26").Value = 20
$ws.Range("L26").Value = 50
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 62
$ws.Range("K27").Value = -6.451612903225
$ws.Range("L27").Value = 34.883720930232
$ws.Range("I28").Value = 5
$ws.Range("K28").Value = 66.666666666666
$ws.Range("L28").Value = -28.571428571428
$ws.Range("M28").Value = 400
$ws.Range("N28").Value = -64.285714285714
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = 66.666666666666
$ws.Range("L29").Value = -16.666666666666
$ws.Range("M29").Value = 400
$ws.Range("N29").Value = -64.285714285714

# --- Cells flipping from numeric to text placeholder ("0" / "***.*") ---
# Force text type with leading apostrophe, then copy number-format/style from a
# cell that already carries the correct "text placeholder" style (s=14).
$ws.Range("G15").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("C18").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("G26").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("H26").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("H26").PasteSpecial(-4122)

# --- Cells flipping from text placeholder to numeric ---
# Assign the numeric value, then copy number-format/style from a cell that
# already carries the correct numeric style (s=16).
$ws.Range("C26").Value = 1
$ws.Range("I15").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("F26").Value = 1
$ws.Range("I15").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("I15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("F28").Value = 1
$ws.Range("I15").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("C29").Value = 1
$ws.Range("I15").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1
$ws.Range("I15").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$excel.CutCopyMode = $false
